$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.571.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.996.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.599'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.14'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.371'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0747'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0988'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.294.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.758'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.014.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.523.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0805'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '222.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.64%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.127'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.92%  '
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.63'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.98%  '
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0942'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.456.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.23%  '
$ws.Range("E44").Value = '  -3.81%  '
$ws.Range("E45").Value = '  -8.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.993'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.53%  '
